$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 10600
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 10600
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 31800
$ws.Range("N48").Value = -32384

$ws.Range("H56").Value = 10600
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 10600
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 31800
$ws.Range("N56").Value = -32868

$ws.Range("H92").Value = 10642.333
$ws.Range("I92").Value = 16986.857
$ws.Range("J92").Value = 1760
$ws.Range("K92").Value = 16986.857
$ws.Range("L92").Value = 1760
$ws.Range("M92").Value = -15738.857
$ws.Range("N92").Value = -4256

$ws.Range("H100").Value = 3408.4546
$ws.Range("I100").Value = 3459.375
$ws.Range("J100").Value = 3379.3572
$ws.Range("K100").Value = 3459.375
$ws.Range("L100").Value = 3379.3572
$ws.Range("M100").Value = -2918.375
$ws.Range("N100").Value = -4461.3572

$ws.Range("H113").Value = 3001.2307
$ws.Range("I113").Value = 2435
$ws.Range("J113").Value = 3486.5715
$ws.Range("K113").Value = 2435
$ws.Range("L113").Value = 3486.5715
$ws.Range("M113").Value = 819
$ws.Range("N113").Value = -9994.5715

$ws.Range("H116").Value = 6898924.5
$ws.Range("I116").Value = 13335201
$ws.Range("J116").Value = 2914.2144
$ws.Range("K116").Value = 13335201
$ws.Range("L116").Value = 2914.2144
$ws.Range("M116").Value = -13331759
$ws.Range("N116").Value = -9798.214400000001

$ws.Range("H127").Value = 1051.6625
$ws.Range("I127").Value = 486
$ws.Range("J127").Value = 1097.527
$ws.Range("K127").Value = 1458
$ws.Range("L127").Value = 3292.581
$ws.Range("M127").Value = 3502
$ws.Range("N127").Value = -13212.581

$ws.Range("H132").Value = 3174.919
$ws.Range("I132").Value = 2839.6
$ws.Range("J132").Value = 3569.4119
$ws.Range("K132").Value = 8518.799999999999
$ws.Range("L132").Value = 10708.2357
$ws.Range("M132").Value = -5988.799999999999

$ws.Range("H141").Value = 3458.132
$ws.Range("I141").Value = 1599.6123
$ws.Range("J141").Value = 26225
$ws.Range("K141").Value = 4798.8369
$ws.Range("L141").Value = 78675
$ws.Range("M141").Value = 381.1630999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 12463
$ws.Range("I36").Value = 4926
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 4926
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = -4580
$ws.Range("N36").Value = -20692

$ws.Range("H44").Value = 34000
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 34000
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 34000
$ws.Range("N44").Value = -34976

$ws.Range("H61").Value = 3572.889
$ws.Range("I61").Value = 2222.2856
$ws.Range("J61").Value = 8300
$ws.Range("K61").Value = 2222.2856
$ws.Range("L61").Value = 8300
$ws.Range("M61").Value = -2010.2856
$ws.Range("N61").Value = -8724

$ws.Range("H80").Value = 30664.834
$ws.Range("I80").Value = 20000
$ws.Range("J80").Value = 35997.25
$ws.Range("K80").Value = 20000
$ws.Range("L80").Value = 35997.25
$ws.Range("M80").Value = -19002
$ws.Range("N80").Value = -37993.25

$ws.Range("H83").Value = 30664.834
$ws.Range("I83").Value = 20000
$ws.Range("J83").Value = 35997.25
$ws.Range("K83").Value = 60000
$ws.Range("L83").Value = 107991.75
$ws.Range("M83").Value = -55008
$ws.Range("N83").Value = -117975.75

$ws.Range("H97").Value = 1616.375
$ws.Range("I97").Value = 1417.1428
$ws.Range("J97").Value = 3011
$ws.Range("K97").Value = 1417.1428
$ws.Range("L97").Value = 3011
$ws.Range("M97").Value = -921.1428000000001
$ws.Range("N97").Value = -4003

$ws.Range("H123").Value = 24302.777
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 24302.777
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 24302.777
$ws.Range("N123").Value = -34102.777

$ws.Range("H132").Value = 2017.0892
$ws.Range("I132").Value = 1696.575
$ws.Range("J132").Value = 2818.375
$ws.Range("K132").Value = 5089.725
$ws.Range("L132").Value = 8455.125
$ws.Range("M132").Value = -2559.725
$ws.Range("N132").Value = -13515.125

$ws.Range("H136").Value = 3572.889
$ws.Range("I136").Value = 2222.2856
$ws.Range("J136").Value = 8300
$ws.Range("K136").Value = 6666.8568
$ws.Range("L136").Value = 24900
$ws.Range("M136").Value = -4116.8568
$ws.Range("N136").Value = -30000

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2634.75
$ws.Range("I99").Value = 1600
$ws.Range("J99").Value = 3255.6
$ws.Range("K99").Value = 1600
$ws.Range("L99").Value = 3255.6
$ws.Range("M99").Value = -102
$ws.Range("N99").Value = -6251.6

$ws.Range("H134").Value = 2593.8538
$ws.Range("I134").Value = 2363.923
$ws.Range("J134").Value = 2992.4
$ws.Range("K134").Value = 7091.768999999999
$ws.Range("L134").Value = 8977.200000000001
$ws.Range("M134").Value = -4556.768999999999
$ws.Range("N134").Value = -14047.2

$ws.Range("H139").Value = 103509.336
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 103509.336
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 103509.336
$ws.Range("N139").Value = -113789.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 102291
$ws.Range("I62").Value = 126863.75
$ws.Range("J62").Value = 4000
$ws.Range("K62").Value = 126863.75
$ws.Range("L62").Value = 4000
$ws.Range("M62").Value = -126239.75
$ws.Range("N62").Value = -5248

$ws.Range("H65").Value = 102291
$ws.Range("I65").Value = 126863.75
$ws.Range("J65").Value = 4000
$ws.Range("K65").Value = 634318.75
$ws.Range("L65").Value = 20000
$ws.Range("M65").Value = -631198.75
$ws.Range("N65").Value = -26240

$ws.Range("H99").Value = 3350.0833
$ws.Range("I99").Value = 3289
$ws.Range("J99").Value = 3533.3333
$ws.Range("K99").Value = 3289
$ws.Range("L99").Value = 3533.3333
$ws.Range("M99").Value = -1791
$ws.Range("N99").Value = -6529.3333

$ws.Range("H126").Value = 3350.0833
$ws.Range("I126").Value = 3289
$ws.Range("J126").Value = 3533.3333
$ws.Range("K126").Value = 9867
$ws.Range("L126").Value = 10599.9999
$ws.Range("M126").Value = -7397
$ws.Range("N126").Value = -15539.9999

$ws.Range("H132").Value = 437532.44
$ws.Range("I132").Value = 677083.9
$ws.Range("J132").Value = 1984.4546
$ws.Range("K132").Value = 2031251.7
$ws.Range("L132").Value = 5953.3638
$ws.Range("M132").Value = -2028721.7

$ws.Range("H134").Value = 1434.5454
$ws.Range("I134").Value = 1181.8529
$ws.Range("J134").Value = 2293.7
$ws.Range("K134").Value = 3545.5587
$ws.Range("L134").Value = 6881.099999999999
$ws.Range("M134").Value = -1010.5587
$ws.Range("N134").Value = -11951.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 300.54544
$ws.Range("I23").Value = 10
$ws.Range("J23").Value = 365.1111
$ws.Range("K23").Value = 30
$ws.Range("L23").Value = 1095.3333
$ws.Range("M23").Value = 205
$ws.Range("N23").Value = -1565.3333

$ws.Range("H134").Value = 5171.926
$ws.Range("I134").Value = 3508.75
$ws.Range("J134").Value = 5872.2104
$ws.Range("K134").Value = 10526.25
$ws.Range("L134").Value = 17616.6312
$ws.Range("M134").Value = -5456.25
$ws.Range("N134").Value = -27756.6312

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 26068.182
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 26068.182
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 26068.182
$ws.Range("N93").Value = -29812.182

$ws.Range("H97").Value = 168525.44
$ws.Range("I97").Value = 73684.14
$ws.Range("J97").Value = 500470
$ws.Range("K97").Value = 73684.14
$ws.Range("L97").Value = 500470
$ws.Range("M97").Value = -73188.14
$ws.Range("N97").Value = -501462

$ws.Range("H109").Value = 11984.857
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 11984.857
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 11984.857
$ws.Range("N109").Value = -14064.857

$ws.Range("H122").Value = 2710.5173
$ws.Range("I122").Value = 2297.35
$ws.Range("J122").Value = 3628.6667
$ws.Range("K122").Value = 6892.049999999999
$ws.Range("L122").Value = 10886.0001
$ws.Range("M122").Value = -4442.049999999999
$ws.Range("N122").Value = -15786.0001

$ws.Range("H132").Value = 1865.6487
$ws.Range("I132").Value = 1248.8889
$ws.Range("J132").Value = 3530.9
$ws.Range("K132").Value = 3746.6667
$ws.Range("L132").Value = 10592.7
$ws.Range("M132").Value = -1216.6667
$ws.Range("N132").Value = -15652.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3683.682
$ws.Range("I40").Value = 4143.6
$ws.Range("J40").Value = 3300.4167
$ws.Range("K40").Value = 4143.6
$ws.Range("L40").Value = 3300.4167
$ws.Range("M40").Value = -4007.6
$ws.Range("N40").Value = -3572.4167

$ws.Range("H61").Value = 20054.75
$ws.Range("I61").Value = 21559.273
$ws.Range("J61").Value = 3505
$ws.Range("K61").Value = 21559.273
$ws.Range("L61").Value = 3505
$ws.Range("M61").Value = -21357.273
$ws.Range("N61").Value = -3909

$ws.Range("H93").Value = 1500
$ws.Range("I93").Value = 1500
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1500
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -252
$ws.Range("N93").ClearContents()

$ws.Range("H100").Value = 2975
$ws.Range("I100").Value = 1975
$ws.Range("J100").Value = 3975
$ws.Range("K100").Value = 1975
$ws.Range("L100").Value = 3975
$ws.Range("M100").Value = -1434
$ws.Range("N100").Value = -5057

$ws.Range("H113").Value = 20054.75
$ws.Range("I113").Value = 21559.273
$ws.Range("J113").Value = 3505
$ws.Range("K113").Value = 21559.273
$ws.Range("L113").Value = 3505
$ws.Range("M113").Value = -19389.273
$ws.Range("N113").Value = -7845

$ws.Range("H122").Value = 12862274
$ws.Range("I122").Value = 14710412
$ws.Range("J122").Value = 11116810
$ws.Range("K122").Value = 44131236
$ws.Range("L122").Value = 33350430
$ws.Range("M122").Value = -44128786
$ws.Range("N122").Value = -33355330

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 36792.332
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 36792.332
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 36792.332
$ws.Range("N16").Value = -37376.332

$ws.Range("H107").Value = 536.38464
$ws.Range("I107").Value = 563.7222
$ws.Range("J107").Value = 474.875
$ws.Range("K107").Value = 1691.1666
$ws.Range("L107").Value = 1424.625
$ws.Range("M107").Value = 228.8334
$ws.Range("N107").Value = -5264.625

$ws.Range("H109").Value = 32344.25
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 32344.25
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 32344.25
$ws.Range("N109").Value = -35118.25

$ws.Range("H112").Value = 500000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 500000
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 500000
$ws.Range("N112").Value = -502954

$ws.Range("H120").Value = 26306.334
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 26306.334
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 26306.334
$ws.Range("N120").Value = -35982.334

$ws.Range("H123").Value = 22150.908
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 22150.908
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 22150.908
$ws.Range("N123").Value = -31950.908
